$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Jean%Dolbeault%xref no email%1,        Gabriel%Turinici%xref no email%1]'
$ws.Range("E4").Value = '[ Nicolas%Hoertel%null%1,      Martin%Blachier%null%1,      Carlos%Blanco%null%1,      Mark%Olfson%null%1,      Marc%Massetti%null%1,      Marina Sánchez%Rico%null%1,      Frédéric%Limosin%null%1,      Henri%Leleu%null%1]'
$ws.Range("E5").Value = '[Tuo%Ji%NULL%1,             Hai-Lian%Chen%NULL%1,             Jing%Xu%NULL%1,             Ling-Ning%Wu%NULL%1,             Jie-Jia%Li%NULL%1,             Kai%Chen%NULL%1,             Gang%Qin%tonygqin@ntu.edu.cn%1]'
$ws.Range("E6").Value = '[Can%Hou%xref no email%1,        Jiaxin%Chen%xref no email%1,        Yaqing%Zhou%xref no email%1,        Lei%Hua%xref no email%1,        Jinxia%Yuan%xref no email%1,        Shu%He%xref no email%1,        Yi%Guo%xref no email%1,        Sheng%Zhang%xref no email%2,        Qiaowei%Jia%xref no email%1,        Chenhui%Zhao%xref no email%1,        Jing%Zhang%xref no email%1,        Guangxu%Xu%xref no email%1,        Enzhi%Jia%xref no email%1]'
$ws.Range("E7").Value = '[Hien%Lau%hlau2@uci.edu%1,             Veria%Khosrawipour%veriakhosrawipour@yahoo.de%1,             Piotr%Kocbach%piotr.kocbach@uwm.edu.pl%1,             Agata%Mikolajczyk%agata.mikolajczyk@upwr.edu.pl%1,             Justyna%Schubert%justyna.schubert@upwr.edu.pl%1,             Jacek%Bania%jacek.bania@upwr.edu.pl%1,             Tanja%Khosrawipour%tkhosrawipour@gmail.com%1]'
$ws.Range("E8").Value = '[Carlo%Signorelli%NULL%1,             Thea%Scognamiglio%NULL%1,             Anna%Odone%NULL%1]'
$ws.Range("E9").Value = '[Jussi%Taipale%xref no email%1,        Paul%Romer%xref no email%1,        Sten%Linnarsson%xref no email%1]'
$ws.Range("E11").Value = '[Nadya%Johanna%NULL%1,             Henrico%Citrawijaya%NULL%1,             Grace%Wangge%NULL%1]'
$ws.Range("E12").Value = '[Gerard J.%Tellis%xref no email%1,        Ashish%Sood%xref no email%1,        Nitish%Sood%xref no email%1]'
$ws.Range("E13").Value = '[Biao%Tang%NULL%1,             Fan%Xia%NULL%1,             Sanyi%Tang%NULL%1,             Nicola Luigi%Bragazzi%NULL%1,             Qian%Li%NULL%1,             Xiaodan%Sun%NULL%1,             Juhua%Liang%NULL%1,             Yanni%Xiao%yxiao@mail.xjtu.edu.cn%1,             Jianhong%Wu%wujh@yorku.ca%1]'
